# New crime data collected - update weekly CompStat report (84th Precinct)
# This script applies the week-over-week edits described in the commit diff:
#  - bumps the "Volume/Number" header and the reporting week dates
#  - updates crime-complaint statistics for rows 15-30 (Week to Date, 28 Day,
#    Year to Date, % change columns), including a few cells whose type flips
#    between a numeric value and a text placeholder ("0" / "***.*")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

# Set a plain numeric value into a cell, keeping/forcing a numeric style by
# copying the number format from a cell that is already styled as a number.
function Set-Num($ws, $ref, $val) {
    $ws.Range($ref).Value2 = $val
}

# Convert a cell to a numeric cell, borrowing the style (cellXf) from a
# known-good numeric cell elsewhere on the sheet so the resulting style index
# matches what a normal numeric cell in that column looks like.
function Set-NumTyped($ws, $ref, $val, $styleSrcRef) {
    $dst = $ws.Range($ref)
    $src = $ws.Range($styleSrcRef)
    $dst.Value2 = $val
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# Convert a cell to a text placeholder cell (shared string "0" or "***.*"),
# borrowing the style (cellXf) from a known-good text placeholder cell.
function Set-TextTyped($ws, $ref, $text, $styleSrcRef) {
    $dst = $ws.Range($ref)
    $src = $ws.Range($styleSrcRef)
    $dst.Value2 = $text
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4163) | Out-Null   # xlPasteValues (force text type)
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Header text updates (Volume/Number + reporting week date range)
# ---------------------------------------------------------------------------

$hdr = $ws.Range("A8")
$hdrChars = $hdr.Characters(21, 2)
$hdrChars.Text = "46"

$week = $ws.Range("C9")
# Edit the later substring first so the earlier substring's character offset
# is not shifted by the date strings changing length (9 -> 10 chars).
$weekChars2 = $week.Characters(47, 10)
$weekChars2.Text = "11/19/2023"
$weekChars1 = $week.Characters(27, 9)
$weekChars1.Text = "11/13/2023"

# ---------------------------------------------------------------------------
# Row 15 (Rape)
# ---------------------------------------------------------------------------
Set-Num $ws "F15" 1
Set-Num $ws "N15" -73.684210526315

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
Set-Num $ws "C16" 5
Set-Num $ws "D16" 1
Set-Num $ws "E16" 400
Set-Num $ws "F16" 11
Set-Num $ws "G16" 6
Set-Num $ws "H16" 83.333333333333
Set-Num $ws "I16" 141
Set-Num $ws "J16" 124
Set-Num $ws "K16" 13.709677419354
Set-Num $ws "L16" 16.528925619834
Set-Num $ws "M16" -11.875
Set-Num $ws "N16" -87.844827586206

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
Set-Num $ws "C17" 4
Set-TextTyped $ws "D17" "0" "C22"
Set-TextTyped $ws "E17" "***.*" "N22"
Set-Num $ws "G17" 17
Set-Num $ws "H17" 23.529411764705
Set-Num $ws "I17" 234
Set-Num $ws "K17" 48.101265822784
Set-Num $ws "L17" 51.948051948051
Set-Num $ws "M17" 127.184466019417
Set-Num $ws "N17" -40.609137055837

# ---------------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------------
Set-Num $ws "C18" 3
Set-TextTyped $ws "D18" "0" "C22"
Set-TextTyped $ws "E18" "***.*" "N22"
Set-Num $ws "G18" 5
Set-Num $ws "H18" 40
Set-Num $ws "I18" 169
Set-Num $ws "K18" 2.424242424242
Set-Num $ws "L18" 29.007633587786
Set-Num $ws "M18" 57.943925233644
Set-Num $ws "N18" -77.006802721088

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
Set-Num $ws "C19" 11
Set-Num $ws "D19" 13
Set-Num $ws "E19" -15.384615384615
Set-Num $ws "F19" 41
Set-Num $ws "H19" -18
Set-Num $ws "I19" 607
Set-Num $ws "J19" 568
Set-Num $ws "K19" 6.866197183098
Set-Num $ws "L19" 29.978586723768
Set-Num $ws "M19" 33.406593406593
Set-Num $ws "N19" -39.117352056168

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
Set-Num $ws "D20" 3
Set-Num $ws "E20" -33.333333333333
Set-Num $ws "G20" 7
Set-Num $ws "H20" -14.285714285714
Set-Num $ws "I20" 71
Set-Num $ws "J20" 59
Set-Num $ws "K20" 20.338983050847
Set-Num $ws "L20" 97.222222222222
Set-Num $ws "M20" 42
Set-Num $ws "N20" -89.339339339339

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
Set-Num $ws "C21" 25
Set-Num $ws "D21" 17
Set-Num $ws "E21" 47.058823529411
Set-Num $ws "F21" 87
Set-Num $ws "G21" 85
Set-Num $ws "H21" 2.352941176470
Set-Num $ws "I21" 1228
Set-Num $ws "J21" 1085
Set-Num $ws "K21" 13.179723502304
Set-Num $ws "L21" 33.623503808487
Set-Num $ws "M21" 39.704209328782
Set-Num $ws "N21" -69.161225514816

# ---------------------------------------------------------------------------
# Row 22 (Transit)
# ---------------------------------------------------------------------------
Set-Num $ws "G22" 2
Set-Num $ws "H22" 0
Set-Num $ws "L22" -12
Set-Num $ws "M22" -32.307692307692

# ---------------------------------------------------------------------------
# Row 23 (Housing)
# ---------------------------------------------------------------------------
Set-TextTyped $ws "D23" "0" "C22"
Set-TextTyped $ws "E23" "***.*" "N22"
Set-Num $ws "L23" -8.695652173913

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
Set-Num $ws "C24" 44
Set-Num $ws "D24" 31
Set-Num $ws "E24" 41.935483870967
Set-Num $ws "F24" 138
Set-Num $ws "G24" 161
Set-Num $ws "H24" -14.285714285714
Set-Num $ws "I24" 1838
Set-Num $ws "J24" 1732
Set-Num $ws "K24" 6.120092378752
Set-Num $ws "L24" 54.713804713804
Set-Num $ws "M24" 34.553440702781

# ---------------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------------
Set-Num $ws "C25" 14
Set-Num $ws "D25" 9
Set-Num $ws "E25" 55.555555555555
Set-Num $ws "F25" 44
Set-Num $ws "G25" 35
Set-Num $ws "H25" 25.714285714285
Set-Num $ws "I25" 378
Set-Num $ws "J25" 329
Set-Num $ws "K25" 14.893617021276
Set-Num $ws "L25" 58.158995815899
Set-Num $ws "M25" 16.666666666666

# ---------------------------------------------------------------------------
# Row 26 (UCR Rape*)
# ---------------------------------------------------------------------------
Set-NumTyped $ws "C26" 1 "F15"
Set-Num $ws "F26" 2
Set-Num $ws "I26" 11
Set-Num $ws "K26" -31.25
Set-Num $ws "L26" -47.619047619047

# ---------------------------------------------------------------------------
# Row 27 (Other Sex Crimes)
# ---------------------------------------------------------------------------
Set-Num $ws "F27" 3
Set-Num $ws "G27" 6
Set-Num $ws "H27" -50
Set-Num $ws "J27" 60
Set-Num $ws "K27" -10

# ---------------------------------------------------------------------------
# Row 28 (Shooting Vic.)
# ---------------------------------------------------------------------------
Set-TextTyped $ws "D28" "0" "C22"
Set-TextTyped $ws "E28" "***.*" "N22"

# ---------------------------------------------------------------------------
# Row 29 (Shooting Inc.)
# ---------------------------------------------------------------------------
Set-TextTyped $ws "D29" "0" "C22"
Set-TextTyped $ws "E29" "***.*" "N22"

# ---------------------------------------------------------------------------
# Row 30 (Hate Crimes)
# ---------------------------------------------------------------------------
Set-NumTyped $ws "D30" 1 "F15"
Set-NumTyped $ws "E30" -100 "K22"
Set-Num $ws "G30" 2
Set-Num $ws "H30" 150
Set-Num $ws "J30" 12
Set-Num $ws "K30" 33.333333333333
Set-Num $ws "L30" 128.571428571429
